# Naming convension duzeltildi, bazi hatalar giderildi.
# Renames "Place_Description"/"Place_Type" table fields to
# "Adress_Description"/"Adress_Type", fixes several naming/column mistakes
# in the "Comment" and "Place" tables, removes a couple of stray rows, and
# tweaks a few cosmetic sheet-view properties (column widths, selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (table headers): Place_Description -> Adress_Description, Place_Type -> Adress_Type
$ws.Range("F1").Value = "Adress_Description"
$ws.Range("G1").Value = "Adress_Type"

# --- Comment table (column C) field renames / reorder
$ws.Range("C2").Value = "user_id  PK    int(10)"
$ws.Range("C3").Value = "place_id   PK    int(10)"
$ws.Range("C4").Value = "text                    nvarchar(10000)"
$ws.Range("C5").Value = "createdOn         datetime"
$ws.Range("C6").Value = "IsActive           bit"
$ws.Range("C7").ClearContents()

# --- Adress_Description / Adress_Type table (columns F, G)
$ws.Range("F2").Value = "adress_desc_id int(10) PK"
$ws.Range("G2").Value = "adress_type_id int(2) PK"
$ws.Range("F3").Value = "adress_name    nvarchar(50)"
$ws.Range("G3").Value = "adress_type  nvarchar(10)"
$ws.Range("F4").Value = "adress_type_id  int(2)  FK"

# --- Place table (column B): drop "mainScore" row, shift remaining rows up
$ws.Range("B4").Value = "adress_desc_id    int(10)  (FK - Adress_Description)"
$ws.Range("B6").Value = "placeImage        ByteArray"
$ws.Range("B7").Value = "createdOn     datetime"
$ws.Range("B8").Value = "IsActive           bit"
$ws.Range("B9").ClearContents()

# --- Column widths
$ws.Columns.Item(2).ColumnWidth = 43.166666666666664
$ws.Columns.Item(3).ColumnWidth = 30.666666666666668

# --- Selection moved from B11 to C3
[void]$ws.Range("C3").Select()
